$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1768.8462
$ws.Range("J17").Value = 1999.8
$ws.Range("L17").Value = 5999.4
$ws.Range("N17").Value = -6335.4
$ws.Range("H33").Value = 308.26666
$ws.Range("I33").Value = 123.083336
$ws.Range("K33").Value = 123.083336
$ws.Range("M33").Value = 105.916664
$ws.Range("H88").Value = 1973.375
$ws.Range("I88").Value = 1193
$ws.Range("J88").Value = 2084.8572
$ws.Range("K88").Value = 1193
$ws.Range("L88").Value = 2084.8572
$ws.Range("M88").Value = -787
$ws.Range("N88").Value = -2896.8572
$ws.Range("H91").Value = 1973.375
$ws.Range("I91").Value = 1193
$ws.Range("J91").Value = 2084.8572
$ws.Range("K91").Value = 1193
$ws.Range("L91").Value = 2084.8572
$ws.Range("M91").Value = 211
$ws.Range("N91").Value = -4892.8572
$ws.Range("H96").Value = 2124.2856
$ws.Range("I96").Value = 1268.2
$ws.Range("J96").Value = 4264.5
$ws.Range("K96").Value = 3804.6
$ws.Range("L96").Value = 12793.5
$ws.Range("M96").Value = -2431.6
$ws.Range("N96").Value = -15539.5
$ws.Range("H97").Value = 6042.6665
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = $null
$ws.Range("H101").Value = 372.2
$ws.Range("I101").Value = 372.2
$ws.Range("K101").Value = 1116.6
$ws.Range("M101").Value = 505.4000000000001
$ws.Range("H116").Value = 8000
$ws.Range("J116").Value = 8000
$ws.Range("L116").Value = 8000
$ws.Range("N116").Value = -14884
$ws.Range("H125").Value = 784.7273
$ws.Range("I125").Value = 894.7778
$ws.Range("J125").Value = 289.5
$ws.Range("K125").Value = 8053.000199999999
$ws.Range("L125").Value = 2605.5
$ws.Range("M125").Value = -5593.000199999999
$ws.Range("N125").Value = -7525.5
$ws.Range("H138").Value = 5221.852
$ws.Range("I138").Value = 2045.25
$ws.Range("K138").Value = 6135.75
$ws.Range("M138").Value = -995.75
$ws.Range("H141").Value = 3259.4
$ws.Range("I141").Value = 2324.25
$ws.Range("K141").Value = 6972.75
$ws.Range("M141").Value = -1792.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 653.2857
$ws.Range("J80").Value = 622.5
$ws.Range("L80").Value = 622.5
$ws.Range("N80").Value = -2618.5
$ws.Range("H83").Value = 653.2857
$ws.Range("J83").Value = 622.5
$ws.Range("L83").Value = 3112.5
$ws.Range("N83").Value = -13096.5
$ws.Range("H86").Value = 2276.5
$ws.Range("I86").Value = 1150.8334
$ws.Range("J86").Value = 5653.5
$ws.Range("K86").Value = 1150.8334
$ws.Range("L86").Value = 5653.5
$ws.Range("M86").Value = -27.83339999999998
$ws.Range("N86").Value = -7899.5
$ws.Range("H89").Value = 2276.5
$ws.Range("I89").Value = 1150.8334
$ws.Range("J89").Value = 5653.5
$ws.Range("K89").Value = 5754.166999999999
$ws.Range("L89").Value = 28267.5
$ws.Range("M89").Value = -138.1669999999995
$ws.Range("N89").Value = -39499.5
$ws.Range("H94").Value = 717.8333
$ws.Range("I94").Value = 717.8333
$ws.Range("K94").Value = 717.8333
$ws.Range("M94").Value = -266.8333
$ws.Range("H134").Value = 5450
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 5450
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 16350
$ws.Range("M134").Value = $null
$ws.Range("N134").Value = -21420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 54.18182
$ws.Range("I7").Value = 76.57143000000001
$ws.Range("K7").Value = 76.57143000000001
$ws.Range("M7").Value = 36.42856999999999
$ws.Range("H58").Value = 4136.6
$ws.Range("J58").Value = 4820.1113
$ws.Range("L58").Value = 4820.1113
$ws.Range("N58").Value = -5226.1113
$ws.Range("H99").Value = 17296.883
$ws.Range("I99").Value = 15096
$ws.Range("J99").Value = 18837.5
$ws.Range("K99").Value = 15096
$ws.Range("L99").Value = 18837.5
$ws.Range("M99").Value = -13598
$ws.Range("N99").Value = -21833.5
$ws.Range("H109").Value = 16089
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").Value = $null
$ws.Range("H122").Value = 7218.3076
$ws.Range("I122").Value = 6422.3335
$ws.Range("J122").Value = 9009.25
$ws.Range("K122").Value = 19267.0005
$ws.Range("L122").Value = 27027.75
$ws.Range("M122").Value = -16817.0005
$ws.Range("N122").Value = -31927.75
$ws.Range("H126").Value = 17296.883
$ws.Range("I126").Value = 15096
$ws.Range("J126").Value = 18837.5
$ws.Range("K126").Value = 45288
$ws.Range("L126").Value = 56512.5
$ws.Range("M126").Value = -42818
$ws.Range("N126").Value = -61452.5
$ws.Range("H134").Value = 3032
$ws.Range("I134").Value = 1531.7
$ws.Range("K134").Value = 4595.1
$ws.Range("M134").Value = -2060.1
$ws.Range("H136").Value = 4136.6
$ws.Range("J136").Value = 4820.1113
$ws.Range("L136").Value = 14460.3339
$ws.Range("N136").Value = -19560.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 108.14286
$ws.Range("I11").Value = 108.14286
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 324.42858
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -184.42858
$ws.Range("N11").Value = $null
$ws.Range("H12").Value = 15.375
$ws.Range("J12").Value = 18.5
$ws.Range("L12").Value = 55.5
$ws.Range("N12").Value = -401.5
$ws.Range("H13").Value = 2712
$ws.Range("I13").Value = 282.66666
$ws.Range("K13").Value = 847.9999799999999
$ws.Range("M13").Value = -679.9999799999999
$ws.Range("H23").Value = 203.5
$ws.Range("I23").Value = 600
$ws.Range("K23").Value = 1800
$ws.Range("M23").Value = -1565
$ws.Range("H26").Value = 416.36365
$ws.Range("I26").Value = 128.25
$ws.Range("K26").Value = 384.75
$ws.Range("M26").Value = -96.75
$ws.Range("H41").Value = 524.5
$ws.Range("I41").Value = 524.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1573.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1235.5
$ws.Range("N41").Value = $null
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = $null
$ws.Range("H80").Value = 4869
$ws.Range("J80").Value = 4707.7144
$ws.Range("L80").Value = 14123.1432
$ws.Range("N80").Value = -15995.1432
$ws.Range("H81").Value = 1846.4
$ws.Range("J81").Value = 1830.5
$ws.Range("L81").Value = 5491.5
$ws.Range("N81").Value = -7737.5
$ws.Range("H83").Value = 4869
$ws.Range("J83").Value = 4707.7144
$ws.Range("L83").Value = 42369.4296
$ws.Range("N83").Value = -51729.4296
$ws.Range("H84").Value = 1846.4
$ws.Range("J84").Value = 1830.5
$ws.Range("L84").Value = 16474.5
$ws.Range("N84").Value = -27706.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2575.818
$ws.Range("I132").Value = 1944.6
$ws.Range("K132").Value = 5833.799999999999
$ws.Range("M132").Value = -3303.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5274.6523
$ws.Range("I16").Value = 3827.3157
$ws.Range("K16").Value = 3827.3157
$ws.Range("M16").Value = -3657.3157
$ws.Range("H22").Value = 470
$ws.Range("I22").Value = 532.5
$ws.Range("K22").Value = 532.5
$ws.Range("M22").Value = -237.5
$ws.Range("H27").Value = 470
$ws.Range("I27").Value = 532.5
$ws.Range("K27").Value = 532.5
$ws.Range("M27").Value = -425.5
$ws.Range("H55").Value = 385.55
$ws.Range("I55").Value = 331.93332
$ws.Range("J55").Value = 546.4
$ws.Range("K55").Value = 331.93332
$ws.Range("L55").Value = 546.4
$ws.Range("M55").Value = -158.93332
$ws.Range("N55").Value = -892.4
$ws.Range("H82").Value = 2928.7
$ws.Range("I82").Value = 3417.4
$ws.Range("J82").Value = 2440
$ws.Range("K82").Value = 3417.4
$ws.Range("L82").Value = 2440
$ws.Range("M82").Value = -3056.4
$ws.Range("N82").Value = -3162
$ws.Range("H85").Value = 2928.7
$ws.Range("I85").Value = 3417.4
$ws.Range("J85").Value = 2440
$ws.Range("K85").Value = 3417.4
$ws.Range("L85").Value = 2440
$ws.Range("M85").Value = -2169.4
$ws.Range("N85").Value = -4936
$ws.Range("H132").Value = 5169.8823
$ws.Range("J132").Value = 5197
$ws.Range("L132").Value = 15591
$ws.Range("N132").Value = -20651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 8595
$ws.Range("I4").Value = 25017.5
$ws.Range("K4").Value = 25017.5
$ws.Range("M4").Value = -24904.5
$ws.Range("H81").Value = 3500
$ws.Range("I81").Value = 3500
$ws.Range("K81").Value = 7000
$ws.Range("M81").Value = -5939
$ws.Range("H84").Value = 3500
$ws.Range("I84").Value = 3500
$ws.Range("K84").Value = 35000
$ws.Range("M84").Value = -29696
$ws.Range("H132").Value = 2153.4
$ws.Range("I132").Value = 2153.4
$ws.Range("K132").Value = 6460.200000000001
$ws.Range("M132").Value = -3930.200000000001

Write-Output "Applied all Seraph_Profits updates"